$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows for the two records (row 2 & row 4) and (row 3 & row 5)
# were swapped in place (same columns A:T, values exchanged between rows).

$rangeA2 = $ws.Range("A2:T2")
$rangeA4 = $ws.Range("A4:T4")
$valuesRow2 = $rangeA2.Value2
$valuesRow4 = $rangeA4.Value2
$rangeA2.Value2 = $valuesRow4
$rangeA4.Value2 = $valuesRow2

$rangeA3 = $ws.Range("A3:T3")
$rangeA5 = $ws.Range("A5:T5")
$valuesRow3 = $rangeA3.Value2
$valuesRow5 = $rangeA5.Value2
$rangeA3.Value2 = $valuesRow5
$rangeA5.Value2 = $valuesRow3
